# Apply the edit described by the commit:
# "Change calc function to return a series of values instead of just final net worth"
#
# This adds a new data series (columns Q:U) to the "Compounded Annually" sheet,
# mirroring the F:I / K:O series but tracking a running "negative net worth"
# scenario with principal growth, and flips which sheet/tab is active.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Compounded Annually")
$ws2 = $wb.Worksheets.Item("Monthly")

# --- New shared string / header -------------------------------------------------
$ws1.Range("Q1").Value = "Annually with Principle +  growth, negative net worth"

# --- Column G width on sheet1 ----------------------------------------------------
$ws1.Columns.Item(7).ColumnWidth = 13.1

# --- Row 2 (seed values) ---------------------------------------------------------
$ws1.Range("Q2").Value = 0
$ws1.Range("R2").Value = -100000
$ws1.Range("S2").Value = 1
$ws1.Range("T2").Value = 40000
$ws1.Range("U2").Value = 1.03

# --- Row 3 ------------------------------------------------------------------------
$ws1.Range("Q3").Value = 1
$ws1.Range("R3").Formula = "=R2*S2+T2"
$ws1.Range("S3").Value = 1
$ws1.Range("T3").Formula = "=T2*U2"
$ws1.Range("U3").Value = 1.03

# --- Row 4 ------------------------------------------------------------------------
$ws1.Range("Q4").Value = 2
$ws1.Range("R4").Formula = "=R3*S3+T3"
$ws1.Range("S4").Value = 1
$ws1.Range("T4").Formula = "=T3*U3"
$ws1.Range("U4").Value = 1.03

# --- Q column: row index 3..58 for rows 5..60 --------------------------------------
for ($r = 5; $r -le 60; $r++) {
    $ws1.Cells.Item($r, 17).Value = ($r - 2)   # column Q = 17
}

# --- S column: 1.07 for rows 5..62 --------------------------------------------------
for ($r = 5; $r -le 62; $r++) {
    $ws1.Cells.Item($r, 19).Value = 1.07       # column S = 19
}

# --- U column: 1.03 constant for rows 2..63 (rows 2-4 already set above) -----------
for ($r = 5; $r -le 63; $r++) {
    $ws1.Cells.Item($r, 21).Value = 1.03       # column U = 21
}

# --- T column: shared growth formula for rows 5..39 --------------------------------
$ws1.Range("T5:T39").Formula = "=T4*U4"

# --- T column: flat literal 40000 from row 40..62 -----------------------------------
for ($r = 40; $r -le 62; $r++) {
    $ws1.Cells.Item($r, 20).Value = 40000      # column T = 20
}

# --- R column: running net-worth formula for rows 5..63 -----------------------------
for ($r = 5; $r -le 63; $r++) {
    $ws1.Cells.Item($r, 18).Formula = "=R" + ($r - 1) + "*S" + ($r - 1) + "+T" + ($r - 1)
}

# --- G column (existing series): convert to a shared formula across G4:G63 ----------
$ws1.Range("G4:G63").Formula = "=G3*H3+I3"

# --- View / selection state ----------------------------------------------------------
[void]$ws2.Activate()
[void]$ws2.Range("K14").Select()

[void]$ws1.Activate()
[void]$ws1.Range("Q2").Select()
